# characteristics_data.xlsx update
# - drop the "ck2.txt" row (the whole row is removed, shifting the ck3.txt
#   row up from row 4 to row 3)
# - bump the "# Vin Values" reading for ck1.txt (C2) from 39 to 70
# - add a new, empty, underlined cell at D5 (info.txt placeholder) and move
#   the active selection there

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ck1.txt's "# Vin Values" reading (C2): 39 -> 70
$ws.Range("C2").Value = 70

# Remove the ck2.txt row entirely (old row 3); ck3.txt row shifts from 4 -> 3
$ws.Rows("3").Delete()

# New placeholder cell D5, empty value, underlined font
$ws.Range("D5").Font.Underline = $true

# Move the selection to the new cell
[void]$ws.Range("D5").Select()

Write-Host "characteristics_data.xlsx updated"
